$d = $word.ActiveDocument

# 1) yth_nama -> kepada_yth
$d.Content.Find.Execute("yth_nama", $true, $false, $false, $false, $false,
                         $true, 1, $false, "kepada_yth", 2)

# 2) The "di " + "yth_lokasi" line ("Kepada Yth. ... di yth_lokasi") needs:
#      "di " (trailing space)  -> "di"  (no trailing space)
#      "yth_lokasi"            -> " lokasi_yth" (leading space)
#    Scope the find/replace to just that paragraph's range so the generic
#    "di " substring elsewhere in the document (e.g. "perkawinan di wilayah")
#    is left untouched.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "di yth_lokasi`r") {
        $rng = $p.Range
        $rng.Find.Execute("di ", $true, $false, $false, $false, $false,
                           $true, 1, $false, "di", 2)
        $rng2 = $p.Range
        $rng2.Find.Execute("yth_lokasi", $true, $false, $false, $false, $false,
                            $true, 1, $false, " lokasi_yth", 2)
    }
}
